$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 6).Value = 1.46
$ws.Cells.Item(2, 9).Value = 14
$ws.Cells.Item(2, 10).Value = 4.1
$ws.Cells.Item(2, 11).Value = 8
$ws.Cells.Item(2, 14).Value = 2.06
$ws.Cells.Item(2, 16).Value = 2.06
$ws.Cells.Item(2, 17).Value = 1.54
$ws.Cells.Item(2, 18).Value = 1.42
$ws.Cells.Item(2, 19).Value = 2.3
$ws.Cells.Item(2, 22).Value = 1.07

# Row 3
$ws.Cells.Item(3, 6).Value = 2.58
$ws.Cells.Item(3, 7).Value = 3.05
$ws.Cells.Item(3, 9).Value = 2.78
$ws.Cells.Item(3, 10).Value = 3.4
$ws.Cells.Item(3, 11).Value = 4.6
$ws.Cells.Item(3, 12).Value = 1.01
$ws.Cells.Item(3, 13).Value = 1.01
$ws.Cells.Item(3, 14).Value = 3.6
$ws.Cells.Item(3, 15).Value = 1.19
$ws.Cells.Item(3, 16).Value = 1.41
$ws.Cells.Item(3, 17).Value = 1.56
$ws.Cells.Item(3, 18).Value = 1.41
$ws.Cells.Item(3, 19).Value = 2.32
$ws.Cells.Item(3, 20).Value = 1.01
$ws.Cells.Item(3, 21).Value = 1.01
$ws.Cells.Item(3, 22).Value = 1.59
$ws.Cells.Item(3, 23).Value = 1.49
$ws.Cells.Item(3, 24).Value = 30
$ws.Cells.Item(3, 25).Value = 20
$ws.Cells.Item(3, 26).Value = 28
$ws.Cells.Item(3, 27).Value = 1000
$ws.Cells.Item(3, 28).Value = 20
$ws.Cells.Item(3, 29).Value = 13
$ws.Cells.Item(3, 30).Value = 17.5
$ws.Cells.Item(3, 31).Value = 36
$ws.Cells.Item(3, 32).Value = 29
$ws.Cells.Item(3, 33).Value = 18
$ws.Cells.Item(3, 34).Value = 21
$ws.Cells.Item(3, 35).Value = 1000
$ws.Cells.Item(3, 36).Value = 1000
$ws.Cells.Item(3, 37).Value = 38
$ws.Cells.Item(3, 38).Value = 1000
$ws.Cells.Item(3, 39).Value = 1000
$ws.Cells.Item(3, 40).Value = 25
$ws.Cells.Item(3, 41).Value = 23

# Row 4
$ws.Cells.Item(4, 6).Value = 1.57
$ws.Cells.Item(4, 8).Value = 6.2
$ws.Cells.Item(4, 9).Value = 8.6
$ws.Cells.Item(4, 11).Value = 1000
$ws.Cells.Item(4, 12).Value = 1.38
$ws.Cells.Item(4, 13).Value = 1.05
$ws.Cells.Item(4, 14).Value = 3.9
$ws.Cells.Item(4, 15).Value = 1.29
$ws.Cells.Item(4, 16).Value = 2
$ws.Cells.Item(4, 17).Value = 1.81
$ws.Cells.Item(4, 18).Value = 1.33
$ws.Cells.Item(4, 19).Value = 2.84
$ws.Cells.Item(4, 20).Value = 1.01
$ws.Cells.Item(4, 21).Value = 1.01
$ws.Cells.Item(4, 22).Value = 1.17
$ws.Cells.Item(4, 23).Value = 2.5
$ws.Cells.Item(4, 24).Value = 1000
$ws.Cells.Item(4, 25).Value = 30
$ws.Cells.Item(4, 26).Value = 70
$ws.Cells.Item(4, 27).Value = 1000
$ws.Cells.Item(4, 28).Value = 12.5
$ws.Cells.Item(4, 29).Value = 14.5
$ws.Cells.Item(4, 30).Value = 38
$ws.Cells.Item(4, 31).Value = 1000
$ws.Cells.Item(4, 32).Value = 11.5
$ws.Cells.Item(4, 33).Value = 12.5
$ws.Cells.Item(4, 34).Value = 32
$ws.Cells.Item(4, 35).Value = 1000
$ws.Cells.Item(4, 36).Value = 22
$ws.Cells.Item(4, 37).Value = 25
$ws.Cells.Item(4, 38).Value = 46
$ws.Cells.Item(4, 39).Value = 1000
$ws.Cells.Item(4, 40).Value = 1000
$ws.Cells.Item(4, 41).Value = 1000

# Row 5
$ws.Cells.Item(5, 6).Value = 4.5
$ws.Cells.Item(5, 7).Value = 6.8
$ws.Cells.Item(5, 8).Value = 1.73
$ws.Cells.Item(5, 9).Value = 2.02
$ws.Cells.Item(5, 10).Value = 3.55
$ws.Cells.Item(5, 11).Value = 5.2
$ws.Cells.Item(5, 14).Value = 1.74
$ws.Cells.Item(5, 15).Value = 1.38
$ws.Cells.Item(5, 16).Value = 1.74
$ws.Cells.Item(5, 17).Value = 1.38
$ws.Cells.Item(5, 19).Value = 3.7
$ws.Cells.Item(5, 20).Value = 1.01
$ws.Cells.Item(5, 21).Value = 1.01
$ws.Cells.Item(5, 22).Value = 1.99
$ws.Cells.Item(5, 23).Value = 1.17
$ws.Cells.Item(5, 24).Value = 1000
$ws.Cells.Item(5, 25).Value = 1000
$ws.Cells.Item(5, 26).Value = 1000
$ws.Cells.Item(5, 27).Value = 1000
$ws.Cells.Item(5, 28).Value = 1000
$ws.Cells.Item(5, 29).Value = 1000
$ws.Cells.Item(5, 30).Value = 1000
$ws.Cells.Item(5, 31).Value = 1000
$ws.Cells.Item(5, 32).Value = 1000
$ws.Cells.Item(5, 33).Value = 1000
$ws.Cells.Item(5, 34).Value = 1000
$ws.Cells.Item(5, 35).Value = 1000
$ws.Cells.Item(5, 36).Value = 1000
$ws.Cells.Item(5, 37).Value = 1000
$ws.Cells.Item(5, 38).Value = 1000
$ws.Cells.Item(5, 39).Value = 1000
$ws.Cells.Item(5, 40).Value = 1000
$ws.Cells.Item(5, 41).Value = 1000

# Row 6
$ws.Cells.Item(6, 6).Value = 1.04
$ws.Cells.Item(6, 7).Value = 1000
$ws.Cells.Item(6, 8).Value = 1.04
$ws.Cells.Item(6, 9).Value = 1000
$ws.Cells.Item(6, 10).Value = 1.01
$ws.Cells.Item(6, 11).Value = 1000
$ws.Cells.Item(6, 16).Value = 1.25
$ws.Cells.Item(6, 17).Value = 1.01

# Row 7
$ws.Cells.Item(7, 6).Value = 1.91
$ws.Cells.Item(7, 9).Value = 5
$ws.Cells.Item(7, 10).Value = 3.3
$ws.Cells.Item(7, 11).Value = 3.7
$ws.Cells.Item(7, 16).Value = 1.73
$ws.Cells.Item(7, 17).Value = 2

# Row 8
$ws.Cells.Item(8, 7).Value = 1.47
$ws.Cells.Item(8, 9).Value = 7.6

# Row 9
$ws.Cells.Item(9, 6).Value = 2.16
$ws.Cells.Item(9, 7).Value = 2.8
$ws.Cells.Item(9, 8).Value = 2.96
$ws.Cells.Item(9, 10).Value = 2.42
$ws.Cells.Item(9, 11).Value = 5.1
$ws.Cells.Item(9, 16).Value = 1.83
$ws.Cells.Item(9, 17).Value = 1.71

# Row 10
$ws.Cells.Item(10, 16).Value = 1.61
$ws.Cells.Item(10, 18).Value = 1.22
$ws.Cells.Item(10, 20).Value = 2.04
$ws.Cells.Item(10, 24).Value = 1000
$ws.Cells.Item(10, 25).Value = 1000
$ws.Cells.Item(10, 26).Value = 980
$ws.Cells.Item(10, 27).Value = 1000
$ws.Cells.Item(10, 28).Value = 1000
$ws.Cells.Item(10, 30).Value = 1000
$ws.Cells.Item(10, 31).Value = 980
$ws.Cells.Item(10, 32).Value = 1000
$ws.Cells.Item(10, 33).Value = 1000
$ws.Cells.Item(10, 34).Value = 1000
$ws.Cells.Item(10, 35).Value = 1000
$ws.Cells.Item(10, 36).Value = 980
$ws.Cells.Item(10, 37).Value = 980
$ws.Cells.Item(10, 38).Value = 1000
$ws.Cells.Item(10, 39).Value = 1000
$ws.Cells.Item(10, 40).Value = 980
$ws.Cells.Item(10, 41).Value = 340

# Row 11
$ws.Cells.Item(11, 6).Value = 4.2
$ws.Cells.Item(11, 7).Value = 8
$ws.Cells.Item(11, 9).Value = 2.12
$ws.Cells.Item(11, 10).Value = 3.25
$ws.Cells.Item(11, 11).Value = 5
$ws.Cells.Item(11, 16).Value = 1.68
$ws.Cells.Item(11, 17).Value = 2.06

# Row 14
$ws.Cells.Item(14, 6).Value = 3.6
$ws.Cells.Item(14, 8).Value = 1.72
$ws.Cells.Item(14, 11).Value = 7.2
$ws.Cells.Item(14, 17).Value = 1.59

# Row 15
$ws.Cells.Item(15, 6).Value = 2.14
$ws.Cells.Item(15, 7).Value = 2.64
$ws.Cells.Item(15, 8).Value = 2.64
$ws.Cells.Item(15, 9).Value = 5.5
$ws.Cells.Item(15, 10).Value = 2.86
$ws.Cells.Item(15, 11).Value = 4.1
$ws.Cells.Item(15, 17).Value = 2.34

# Row 16
$ws.Cells.Item(16, 6).Value = 1.21
$ws.Cells.Item(16, 7).Value = 1000
$ws.Cells.Item(16, 8).Value = 1.21
$ws.Cells.Item(16, 9).Value = 1000
$ws.Cells.Item(16, 10).Value = 1.01
$ws.Cells.Item(16, 11).Value = 1000
$ws.Cells.Item(16, 16).Value = 1.24
$ws.Cells.Item(16, 17).Value = 1.01

# Row 17
$ws.Cells.Item(17, 6).Value = 1.04
$ws.Cells.Item(17, 7).Value = 1000
$ws.Cells.Item(17, 8).Value = 1.04
$ws.Cells.Item(17, 9).Value = 1000
$ws.Cells.Item(17, 10).Value = 1.01
$ws.Cells.Item(17, 11).Value = 1000
$ws.Cells.Item(17, 16).Value = 1.24
$ws.Cells.Item(17, 17).Value = 1.01

# Row 19
$ws.Cells.Item(19, 6).Value = 2.46
$ws.Cells.Item(19, 7).Value = 3.05
$ws.Cells.Item(19, 8).Value = 2.86
$ws.Cells.Item(19, 10).Value = 3
$ws.Cells.Item(19, 16).Value = 1.25
$ws.Cells.Item(19, 17).Value = 1.01

# Row 20
$ws.Cells.Item(20, 8).Value = 2.52
$ws.Cells.Item(20, 10).Value = 2.9
$ws.Cells.Item(20, 16).Value = 1.25
$ws.Cells.Item(20, 17).Value = 1.01

# Row 21
$ws.Cells.Item(21, 14).Value = 2.98
$ws.Cells.Item(21, 16).Value = 1.63
$ws.Cells.Item(21, 18).Value = 1.23
$ws.Cells.Item(21, 20).Value = 2.04
$ws.Cells.Item(21, 24).Value = 980
$ws.Cells.Item(21, 25).Value = 1000
$ws.Cells.Item(21, 26).Value = 1000
$ws.Cells.Item(21, 27).Value = 1000
$ws.Cells.Item(21, 28).Value = 1000
$ws.Cells.Item(21, 30).Value = 1000
$ws.Cells.Item(21, 31).Value = 1000
$ws.Cells.Item(21, 32).Value = 1000
$ws.Cells.Item(21, 33).Value = 980
$ws.Cells.Item(21, 34).Value = 1000
$ws.Cells.Item(21, 35).Value = 1000
$ws.Cells.Item(21, 36).Value = 1000
$ws.Cells.Item(21, 37).Value = 65
$ws.Cells.Item(21, 38).Value = 1000
$ws.Cells.Item(21, 39).Value = 1000
$ws.Cells.Item(21, 40).Value = 65
$ws.Cells.Item(21, 41).Value = 1000

# Row 22
$ws.Cells.Item(22, 6).Value = 1.93
$ws.Cells.Item(22, 7).Value = 1.94
$ws.Cells.Item(22, 14).Value = 6
$ws.Cells.Item(22, 18).Value = 1.71
$ws.Cells.Item(22, 19).Value = 2.36
$ws.Cells.Item(22, 26).Value = 42
$ws.Cells.Item(22, 27).Value = 90
$ws.Cells.Item(22, 31).Value = 44
$ws.Cells.Item(22, 38).Value = 26
$ws.Cells.Item(22, 39).Value = 60
$ws.Cells.Item(22, 40).Value = 9
$ws.Cells.Item(22, 41).Value = 40

# Row 23
$ws.Cells.Item(23, 8).Value = 3.55
$ws.Cells.Item(23, 15).Value = 1.41
$ws.Cells.Item(23, 16).Value = 1.74
$ws.Cells.Item(23, 21).Value = 1.99
$ws.Cells.Item(23, 25).Value = 12
$ws.Cells.Item(23, 27).Value = 80
$ws.Cells.Item(23, 30).Value = 15
$ws.Cells.Item(23, 31).Value = 48
$ws.Cells.Item(23, 34).Value = 20
$ws.Cells.Item(23, 35).Value = 70
$ws.Cells.Item(23, 38).Value = 46
$ws.Cells.Item(23, 40).Value = 25

# Row 24
$ws.Cells.Item(24, 6).Value = 2.72
$ws.Cells.Item(24, 10).Value = 3.3
$ws.Cells.Item(24, 11).Value = 3.6
$ws.Cells.Item(24, 12).Value = 1.01
$ws.Cells.Item(24, 14).Value = 1.04
$ws.Cells.Item(24, 16).Value = 1.31
$ws.Cells.Item(24, 17).Value = 1.32
$ws.Cells.Item(24, 18).Value = 1.31
$ws.Cells.Item(24, 19).Value = 2.84
$ws.Cells.Item(24, 20).Value = 1.01
$ws.Cells.Item(24, 22).Value = 1.54
$ws.Cells.Item(24, 24).Value = 1000
$ws.Cells.Item(24, 25).Value = 1000
$ws.Cells.Item(24, 26).Value = 1000
$ws.Cells.Item(24, 27).Value = 1000
$ws.Cells.Item(24, 28).Value = 1000
$ws.Cells.Item(24, 29).Value = 1000
$ws.Cells.Item(24, 30).Value = 1000
$ws.Cells.Item(24, 31).Value = 1000
$ws.Cells.Item(24, 32).Value = 1000
$ws.Cells.Item(24, 33).Value = 1000
$ws.Cells.Item(24, 34).Value = 1000
$ws.Cells.Item(24, 35).Value = 1000
$ws.Cells.Item(24, 36).Value = 1000
$ws.Cells.Item(24, 37).Value = 1000
$ws.Cells.Item(24, 38).Value = 1000
$ws.Cells.Item(24, 39).Value = 1000
$ws.Cells.Item(24, 40).Value = 1000
$ws.Cells.Item(24, 41).Value = 1000

# Row 25
$ws.Cells.Item(25, 7).Value = 1000
$ws.Cells.Item(25, 8).Value = 1.04
$ws.Cells.Item(25, 9).Value = 1000
$ws.Cells.Item(25, 10).Value = 1.02
$ws.Cells.Item(25, 11).Value = 1000
$ws.Cells.Item(25, 12).Value = 1.01
$ws.Cells.Item(25, 14).Value = 1.25
$ws.Cells.Item(25, 15).Value = 1.08
$ws.Cells.Item(25, 18).Value = 1.13
$ws.Cells.Item(25, 19).Value = 1.44
$ws.Cells.Item(25, 20).Value = 1.01
$ws.Cells.Item(25, 21).Value = 1.01
$ws.Cells.Item(25, 23).Value = 1.62
$ws.Cells.Item(25, 41).Value = 80
